# The single "Fecha y hora" header column is being split into two separate
# columns: "Fecha" and "Hora". Insert a new column before the old B1
# ("Cliente") so the remaining headers (Cliente, Forma de pago, Monto) shift
# right by one, then relabel A1/B1 as the two new headers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").EntireColumn.Insert() | Out-Null

$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Hora"
